# Form the consolidated report: fill in the computed "Absent" values
# for the rows where they were still blank/placeholder, and correct the
# values that were computed incorrectly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
